$wb = $excel.ActiveWorkbook

# --- Summary sheet: remove the two trailing blank formatted rows (6 & 7) ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Rows.Item(6).Delete()
$wsSummary.Rows.Item(6).Delete()
$wsSummary.Range("E4").Select()

# --- Repayment Schedule sheet: insert a new blank separator column (N) ---
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns.Item(14).Insert()

# give the new header cell (N1) the same shaded header style as the rest of row 1
$wsSchedule.Range("A1").Copy()
$wsSchedule.Cells.Item(1, 14).PasteSpecial(-4122)

# match the column width used for the new separator column
$wsSchedule.Columns.Item(14).ColumnWidth = 9.140625

$wsSchedule.Range("M18").Select()

# --- make "Repayment Schedule" the active tab of the workbook ---
$wsSchedule.Activate()
